$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the destination email value first (this creates the shared string for
# the email address before the header labels, matching the expected
# shared-string table order), then turn it into a mailto: hyperlink, which
# also applies the built-in "Hyperlink" cell style to D2.
$ws.Range("D2").Value = "burcea.bogdan.madalin@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:burcea.bogdan.madalin@gmail.com")

# New column headers
$ws.Range("D1").Value = "DestinationEmail"
$ws.Range("E1").Value = "SourceEmailCredential"

# Source email credential (Windows Credential Manager entry name)
$ws.Range("E2").Value = "bot_gmail"

# Extra styled (Hyperlink look) but otherwise empty cell further down the sheet
$ws.Range("D5").Style = "Hyperlink"

# Resize the new columns
$ws.Columns.Item(4).ColumnWidth = 38.666666666666664
$ws.Columns.Item(5).ColumnWidth = 41.166666666666664

# Update the active selection
$ws.Range("E7").Select()
